# Auto-generated Excel COM-interop script to apply market-price refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 338.79166
$ws.Range("I33").Value = 376.11765
$ws.Range("K33").Value = 376.11765
$ws.Range("M33").Value = -147.11765

$ws.Range("H93").Value = 34995
$ws.Range("J93").Value = 34995
$ws.Range("L93").Value = 34995
$ws.Range("N93").Value = -39987

$ws.Range("H112").Value = 29665.053
$ws.Range("J112").Value = 36141.258
$ws.Range("L112").Value = 108423.774
$ws.Range("N112").Value = -110639.774

$ws.Range("H116").Value = 9019.75
$ws.Range("I116").Value = 11651.5
$ws.Range("K116").Value = 11651.5
$ws.Range("M116").Value = -8209.5

$ws.Range("H121").Value = 2156
$ws.Range("J121").Value = 2156
$ws.Range("L121").Value = 6468
$ws.Range("N121").Value = -9962

$ws.Range("H132").Value = 52268.35
$ws.Range("I132").Value = 2090.8
$ws.Range("K132").Value = 6272.400000000001
$ws.Range("M132").Value = -3742.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H61").Value = 5056.9287
$ws.Range("J61").Value = 5366.5
$ws.Range("L61").Value = 5366.5
$ws.Range("N61").Value = -5790.5

$ws.Range("H110").Value = 5970.7036
$ws.Range("I110").Value = 6350.45
$ws.Range("K110").Value = 6350.45
$ws.Range("M110").Value = -4305.45

$ws.Range("H122").Value = 1858
$ws.Range("I122").Value = 1355.8889
$ws.Range("K122").Value = 4067.6667
$ws.Range("M122").Value = -1617.6667

$ws.Range("H136").Value = 5056.9287
$ws.Range("J136").Value = 5366.5
$ws.Range("L136").Value = 16099.5
$ws.Range("N136").Value = -21199.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2060.3
$ws.Range("I94").Value = 1407.7142
$ws.Range("K94").Value = 1407.7142
$ws.Range("M94").Value = -956.7141999999999

$ws.Range("H105").Value = 1638.1364
$ws.Range("I105").Value = 1713.8948
$ws.Range("K105").Value = 1713.8948
$ws.Range("M105").Value = 33.10519999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1332.3158
$ws.Range("I105").Value = 1285.4445
$ws.Range("K105").Value = 1285.4445
$ws.Range("M105").Value = 461.5554999999999

$ws.Range("H141").Value = 494130.2
$ws.Range("I141").Value = 70000
$ws.Range("K141").Value = 70000
$ws.Range("M141").Value = -64820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62.64706
$ws.Range("I2").Value = 64.066666
$ws.Range("J2").Value = 52
$ws.Range("K2").Value = 384.399996
$ws.Range("L2").Value = 312
$ws.Range("M2").Value = -271.399996
$ws.Range("N2").Value = -538

$ws.Range("H33").Value = 244.33333
$ws.Range("I33").Value = 20
$ws.Range("J33").Value = 356.5
$ws.Range("K33").Value = 120
$ws.Range("L33").Value = 2139
$ws.Range("M33").Value = 163
$ws.Range("N33").Value = -2705

$ws.Range("H49").Value = 1157.3334
$ws.Range("I49").Value = 648.3333
$ws.Range("J49").Value = 1666.3334
$ws.Range("K49").Value = 1944.9999
$ws.Range("L49").Value = 4999.0002
$ws.Range("M49").Value = -1788.9999
$ws.Range("N49").Value = -5311.0002

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H82").Value = 35769
$ws.Range("I82").Value = 13456
$ws.Range("J82").Value = 52503.75
$ws.Range("K82").Value = 40368
$ws.Range("L82").Value = 157511.25
$ws.Range("M82").Value = -39962
$ws.Range("N82").Value = -158323.25

$ws.Range("H85").Value = 35769
$ws.Range("I85").Value = 13456
$ws.Range("J85").Value = 52503.75
$ws.Range("K85").Value = 40368
$ws.Range("L85").Value = 157511.25
$ws.Range("M85").Value = -38964
$ws.Range("N85").Value = -160319.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H61").Value = 3781.4814
$ws.Range("I61").Value = 3602.1904
$ws.Range("K61").Value = 3602.1904
$ws.Range("M61").Value = -3400.1904

$ws.Range("H68").Value = 2651.8462
$ws.Range("I68").Value = 2708.25
$ws.Range("J68").Value = 1975
$ws.Range("K68").Value = 2708.25
$ws.Range("L68").Value = 1975
$ws.Range("M68").Value = -1959.25
$ws.Range("N68").Value = -3473

$ws.Range("H71").Value = 2651.8462
$ws.Range("I71").Value = 2708.25
$ws.Range("J71").Value = 1975
$ws.Range("K71").Value = 13541.25
$ws.Range("L71").Value = 9875
$ws.Range("M71").Value = -9797.25
$ws.Range("N71").Value = -17363

$ws.Range("H113").Value = 3781.4814
$ws.Range("I113").Value = 3602.1904
$ws.Range("K113").Value = 3602.1904
$ws.Range("M113").Value = -1432.1904

$ws.Range("H122").Value = 4682.826
$ws.Range("I122").Value = 4336.2144
$ws.Range("K122").Value = 13008.6432
$ws.Range("M122").Value = -10558.6432

$ws.Range("H132").Value = 4706.294
$ws.Range("I132").Value = 4485.643
$ws.Range("K132").Value = 13456.929
$ws.Range("M132").Value = -10926.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2930
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 3757.1428
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 3757.1428
$ws.Range("M4").Value = -887
$ws.Range("N4").Value = -3983.1428

$ws.Range("H62").Value = 18999.334
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 18999.334
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H69").Value = 44940.75
$ws.Range("I69").Value = 43382
$ws.Range("K69").Value = 43382
$ws.Range("M69").Value = -42633

$ws.Range("H72").Value = 44940.75
$ws.Range("I72").Value = 43382
$ws.Range("K72").Value = 130146
$ws.Range("M72").Value = -126402

$ws.Range("H132").Value = 1998.9375
$ws.Range("I132").Value = 1927.3572
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5782.071599999999
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3252.071599999999
$ws.Range("N132").Value = -12560
